# GEPEP calibration for iteration 1 of GFS for GSW Spring lambing
#
# - Low!I2:I31   0.035 -> 0.03
# - Low!K2:K31   9e-8  -> 0   (and drop the now-unused "40% - Accent6 + 0.00E+00"
#                               number-format style, reverting these cells to the
#                               plain "40% - Accent6" / General style)
# - BestBet becomes the active / selected sheet again (was Low)
# - Low's remembered selection moves from K2:K31 to L2

$wb = $excel.ActiveWorkbook

$low = $wb.Worksheets.Item("Low")

$low.Range("I2:I31").Value = 0.03

$low.Range("K2:K31").Style = "40% - Accent6"
$low.Range("K2:K31").Value = 0

$low.Range("L2").Select() | Out-Null

$bestBet = $wb.Worksheets.Item("BestBet")
$bestBet.Activate() | Out-Null
